$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

$ws.Range("G2").Value = 0.01854753494262695
$ws.Range("H2").Value = 0.002236604690551758
$ws.Range("I2").Value = 0.02078413963317871

$ws.Range("G3").Value = 0.0181126594543457
$ws.Range("H3").Value = 0.00431060791015625
$ws.Range("I3").Value = 0.02242326736450195

$ws.Range("G4").Value = 0.01540899276733398
$ws.Range("H4").Value = 0.002544641494750977
$ws.Range("I4").Value = 0.01795363426208496
